# Player.xlsx / "Property" sheet edit
# - Queue-lock related stat rows (MAXHP..BUFF_GATE, rows 44-67) no longer
#   get auto-saved: flip their "Save" column (E) from TRUE to FALSE.
# - Rows 76-77 (GameID / GateID, the newest fields) are no longer
#   "pending/new" so their red-on-yellow highlight formatting is cleared;
#   GameID's Id cell (A76) keeps the Text number format.
# - Move the sheet's active selection to H78 (bottom of the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")
$ws.Activate()

# Turn off "Save" for the queue-lock related properties (rows 44-67).
$ws.Range("E44:E67").Value = $false

# Un-highlight the recently-added GameID / GateID rows.
$ws.Range("A76:J77").ClearFormats()
$ws.Range("A76").NumberFormat = "@"

# Update the saved selection/active cell.
$ws.Range("H78").Select()
